# résultats_CL_68_param_cosmo(à compléter).xlsx
# - new plots of residuals + update OmegaM_w + creation of ppt
#   'Récapitulatif chi2,HD,résidus'
#
# Concretely, on Sheet1:
#  * Fill in the new "Flat wCDM" results block in row 7 (title + fitted
#    H0/Omegam/OmegaLambda-style values) and row 8 (Omegam_w restraining
#    parameter line, with the ABS() comparison formulas against the
#    Brout et al. reference values).
#  * Drop the old red/blue placeholder font from A7/A8 now that they hold
#    real data.
#  * Leave a note on G8 (the Omegam_w restraining-power delta) pointing
#    out how good the fit looks compared to Brout et al.
#  * Leave the selection on I9, matching where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 / Row 8: new "Flat wCDM" block -------------------------------
# Write the two row labels first so the new shared strings land in the
# same order the workbook ends up with.
$ws.Range("A7").Value = "Résultats attendus pour Flat wCDM (Brout et al. 2022)"
$ws.Range("A8").Value = "Omegam_w(pour flat wCDM)"

# Row 7: fitted values for the flat wCDM run.
$ws.Range("B7").Value = 73.5
$ws.Range("C7").Value = 1.1
$ws.Range("D7").Value = 1.1
$ws.Range("E7").Value = 0.309
$ws.Range("F7").Value = 0.063
$ws.Range("G7").Value = -0.069
$ws.Range("H7").Value = "(w=-0,90) 0,691"
$ws.Range("I7").Value = "(0,14 pour w) 0,069"
$ws.Range("J7").Value = "(0,14 pour w) -0,063"

# Row 8: Omegam_w restraining parameter vs. the Brout reference values.
$ws.Range("E8").Value = 0.301
$ws.Range("F8").Formula = "=ABS(E8-0.308)"
$ws.Range("G8").Formula = "=ABS(E8-0.29)"
$ws.Range("H8").Value = -0.92
$ws.Range("I8").Value = "range trop restreint"
$ws.Range("J8").Value = "range trop restreint"

# The A7/A8 cells used to carry a leftover red/blue placeholder font;
# now that they hold real results, clear that back to the normal style.
$ws.Range("A7").ClearFormats()
$ws.Range("A8").ClearFormats()

# --- New review comment on G8 -------------------------------------------
$commentText = "Emile DOSSO:" + [char]10 + [char]10 + "étonamment précis en comparaison aux résultats de Brout…" + [char]10
$ws.Range("G8").AddComment($commentText) | Out-Null

# --- Column widths widened to fit the new text in I/J --------------------
$ws.Columns.Item(9).ColumnWidth = 18.73
$ws.Columns.Item(10).ColumnWidth = 20.45

# --- Leave the selection where the author left it ------------------------
$ws.Range("I9").Select() | Out-Null
